$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Insert()

Write-Host "UsedRange after insert: $($ws.UsedRange.Address())"
